$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.500.82'
$ws.Range('D3').Value = '3.032.47'
$ws.Range('E3').Value = '  -3.94%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '210.14'
$ws.Range('E5').Value = '  -2.76%  '
$ws.Range('D6').Value = '612.55'
$ws.Range('E6').Value = '  -3.54%  '
$ws.Range('D7').Value = '0.359'
$ws.Range('E7').Value = '  -9.83%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.880'
$ws.Range('E8').Value = '  +20.82%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').Value = '3.029.96'
$ws.Range('E10').Value = '  -4.06%  '
$ws.Range('E11').Value = '  +19.03%  '
$ws.Range('D12').Value = '0.187'
$ws.Range('E12').Value = '  +4.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000236'
$ws.Range('D14').Value = '5.34'
$ws.Range('E14').Value = '  +0.45%  '
$ws.Range('D15').Value = '88.543.49'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').Value = '31.83'
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('D17').Value = '3.590.76'
$ws.Range('E17').Value = '  -3.56%  '
$ws.Range('D18').Value = '3.024.17'
$ws.Range('E18').Value = '  -3.82%  '
$ws.Range('D19').Value = '3.34'
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000213'
$ws.Range('E20').Value = '  -7.67%  '
$ws.Range('D21').Value = '13.33'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').Value = '424.71'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').Value = '4.98'
$ws.Range('E23').Value = '  +1.15%  '
$ws.Range('D24').Value = '8.14'
$ws.Range('E24').Value = '  -3.39%  '
$ws.Range('D25').Value = '5.36'
$ws.Range('E25').Value = '  -0.88%  '
$ws.Range('D26').Value = '83.29'
$ws.Range('E26').Value = '  +3.78%  '
$ws.Range('D27').Value = '11.64'
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('D28').Value = '3.187.99'
$ws.Range('E28').Value = '  -3.43%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '0.162'
$ws.Range('E30').Value = '  +0.78%  '
$ws.Range('E31').Value = '  +1.63%  '
$ws.Range('D32').Value = '8.17'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('D33').Value = '501.28'
$ws.Range('E33').Value = '  -2.79%  '
$ws.Range('D34').Value = '3.61'
$ws.Range('E34').Value = '  -10.65%  '
$ws.Range('D35').Value = '6.61'
$ws.Range('E35').Value = '  -7.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.70'
$ws.Range('E36').Value = '  +3.27%  '
$ws.Range('E37').Value = '  -2.86%  '
$ws.Range('E38').Value = '  -7.24%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.130'
$ws.Range('E39').Value = '  -9.62%  '
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').Value = '22.22'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.360'
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('E44').Value = '  +7.98%  '
$ws.Range('D45').Value = '1.82'
$ws.Range('E45').Value = '  -3.17%  '
$ws.Range('D46').Value = '145.66'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '43.26'
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('D48').Value = '0.0673'
$ws.Range('E48').Value = '  +7.76%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '160.06'
$ws.Range('E49').Value = '  -4.00%  '
$ws.Range('D50').Value = '1.21'
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').Value = '4.04'
$ws.Range('E51').Value = '  +2.20%  '
